$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.248.91"
$ws.Range("E2").Value = "  +2.78%  "
$ws.Range("D3").Value = "3.813.07"
$ws.Range("E3").Value = "  +1.55%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("D7").Value = "3.813.67"
$ws.Range("E7").Value = "  +1.62%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  -2.02%  "
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000264"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").Value = "4.451.15"
$ws.Range("E15").Value = "  +1.43%  "
$ws.Range("D16").Value = "3.807.82"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D17").Value = "69.257.19"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.54%  "
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "471.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = "  -1.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.07%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "3.960.25"
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("E31").Value = "  -2.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.15%  "
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.58%  "
$ws.Range("D37").Value = "3.766.59"
$ws.Range("E37").Value = "  +1.19%  "
$ws.Range("E38").Value = "  -1.76%  "
$ws.Range("E39").Value = "  -6.95%  "
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("E42").Value = "  +0.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.310"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +11.07%  "
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "402.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "145.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.81%  "

Write-Output "Updated cryptos list"
